$wb = $excel.ActiveWorkbook

# --- 1. Metadata sheet: bump the "Date" property value ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-09-23T14:10:57+00:00"

# --- 2. Elements sheet: append a new row describing the
#        "quantitePrescrite" element. The new row reuses the same
#        shape/formatting as the previous last row (row 11), so start by
#        cloning that row (formats + values) into the new row 12, then
#        overwrite just the cells that actually differ. ---
$elements = $wb.Worksheets.Item("Elements")

$elements.Range("A11:AJ11").Copy()
$elements.Range("A12:AJ12").PasteSpecial(-4122)

$elements.Range("A11:AJ11").Copy()
$elements.Range("A12:AJ12").PasteSpecial(-4163)

$excel.CutCopyMode = 0

$newRow = 12
$elements.Range("A" + $newRow).Value = "fr-ligne-prescription.quantitePrescrite"
$elements.Range("B" + $newRow).Value = "fr-ligne-prescription.quantitePrescrite"
$elements.Range("K" + $newRow).Value = "Quantity`n"
$elements.Range("L" + $newRow).Value = "Quantité totale de traitement prescrite, indépendamment des séquences. Cela permet d'aider la dispensation dans le cas d'une unité non convertible (exemple : crème avec une posologie en ""application"")"
$elements.Range("M" + $newRow).Value = "Quantité totale de traitement prescrite, indépendamment des séquences. Cela permet d'aider la dispensation dans le cas d'une unité non convertible (exemple : crème avec une posologie en ""application"")"
$elements.Range("AF" + $newRow).Value = "fr-ligne-prescription.quantitePrescrite"
